# Auto-generated edit script applying numeric value updates to the Goblin_Profits workbook
# (market price / profit recompute across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 769456.7
$ws.Range("I2").Value = 833574.7
$ws.Range("J2").Value = 41
$ws.Range("K2").Value = 833574.7
$ws.Range("L2").Value = 41
$ws.Range("M2").Value = -833461.7
$ws.Range("N2").Value = -267
$ws.Range("H11").Value = 60.941177
$ws.Range("I11").Value = 60.941177
$ws.Range("K11").Value = 60.941177
$ws.Range("M11").Value = 79.05882299999999
$ws.Range("H40").Value = 5000
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 5000
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 5000
$ws.Range("M40").Value = ""
$ws.Range("N40").Value = -5350
$ws.Range("H86").Value = 3403.55
$ws.Range("I86").Value = 3010.7646
$ws.Range("J86").Value = 5629.3335
$ws.Range("K86").Value = 3010.7646
$ws.Range("L86").Value = 5629.3335
$ws.Range("M86").Value = -1887.7646
$ws.Range("N86").Value = -7875.3335
$ws.Range("H89").Value = 3403.55
$ws.Range("I89").Value = 3010.7646
$ws.Range("J89").Value = 5629.3335
$ws.Range("K89").Value = 15053.823
$ws.Range("L89").Value = 28146.6675
$ws.Range("M89").Value = -9437.823
$ws.Range("N89").Value = -39378.6675
$ws.Range("H96").Value = 715446.6
$ws.Range("I96").Value = 1111688.2
$ws.Range("J96").Value = 2211.8
$ws.Range("K96").Value = 3335064.6
$ws.Range("L96").Value = 6635.400000000001
$ws.Range("M96").Value = -3333691.6
$ws.Range("N96").Value = -9381.400000000001
$ws.Range("H131").Value = 6437.636
$ws.Range("I131").Value = 755
$ws.Range("K131").Value = 2265
$ws.Range("M131").Value = 2775
$ws.Range("H137").Value = 1907.2593
$ws.Range("I137").Value = 1563.625
$ws.Range("K137").Value = 4690.875
$ws.Range("M137").Value = -2140.875
$ws.Range("H138").Value = 4277.602
$ws.Range("I138").Value = 3993.2424
$ws.Range("J138").Value = 4434
$ws.Range("K138").Value = 11979.7272
$ws.Range("L138").Value = 13302
$ws.Range("M138").Value = -6839.727200000001
$ws.Range("N138").Value = -23582
$ws.Range("H141").Value = 3567.6667
$ws.Range("I141").Value = 3786
$ws.Range("K141").Value = 11358
$ws.Range("M141").Value = -6178

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H29").Value = 4166.6665
$ws.Range("J29").Value = 4166.6665
$ws.Range("L29").Value = 4166.6665
$ws.Range("N29").Value = -4782.6665
$ws.Range("H32").Value = 2505.7837
$ws.Range("I32").Value = 1879.6086
$ws.Range("K32").Value = 1879.6086
$ws.Range("M32").Value = -1592.6086
$ws.Range("H61").Value = 1642.5686
$ws.Range("I61").Value = 985.4878
$ws.Range("J61").Value = 4336.6
$ws.Range("K61").Value = 985.4878
$ws.Range("L61").Value = 4336.6
$ws.Range("M61").Value = -773.4878
$ws.Range("N61").Value = -4760.6
$ws.Range("H74").Value = 2075.1353
$ws.Range("I74").Value = 2022.2858
$ws.Range("J74").Value = 3000
$ws.Range("K74").Value = 2022.2858
$ws.Range("L74").Value = 3000
$ws.Range("M74").Value = -1148.2858
$ws.Range("N74").Value = -4748
$ws.Range("H77").Value = 2075.1353
$ws.Range("I77").Value = 2022.2858
$ws.Range("J77").Value = 3000
$ws.Range("K77").Value = 10111.429
$ws.Range("L77").Value = 15000
$ws.Range("M77").Value = -5743.429
$ws.Range("N77").Value = -23736
$ws.Range("H132").Value = 3238.1724
$ws.Range("I132").Value = 3569.2917
$ws.Range("K132").Value = 10707.8751
$ws.Range("M132").Value = -8177.875100000001
$ws.Range("H136").Value = 1642.5686
$ws.Range("I136").Value = 985.4878
$ws.Range("J136").Value = 4336.6
$ws.Range("K136").Value = 2956.4634
$ws.Range("L136").Value = 13009.8
$ws.Range("M136").Value = -406.4634000000001
$ws.Range("N136").Value = -18109.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 7900.353
$ws.Range("I20").Value = 9407.846
$ws.Range("K20").Value = 9407.846
$ws.Range("M20").Value = -9160.846
$ws.Range("H86").Value = 27785894
$ws.Range("I86").Value = 9698.615
$ws.Range("J86").Value = 100004000
$ws.Range("K86").Value = 9698.615
$ws.Range("L86").Value = 100004000
$ws.Range("M86").Value = -8575.615
$ws.Range("N86").Value = -100006246
$ws.Range("H89").Value = 27785894
$ws.Range("I89").Value = 9698.615
$ws.Range("J89").Value = 100004000
$ws.Range("K89").Value = 48493.075
$ws.Range("L89").Value = 500020000
$ws.Range("M89").Value = -42877.075
$ws.Range("N89").Value = -500031232
$ws.Range("H134").Value = 1725.6052
$ws.Range("I134").Value = 1424.1875
$ws.Range("J134").Value = 3333.1667
$ws.Range("K134").Value = 4272.5625
$ws.Range("L134").Value = 9999.500100000001
$ws.Range("M134").Value = -1737.5625
$ws.Range("N134").Value = -15069.5001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 26666
$ws.Range("J9").Value = 26666
$ws.Range("L9").Value = 26666
$ws.Range("N9").Value = -27002
$ws.Range("H58").Value = 1305.88
$ws.Range("I58").Value = 1296
$ws.Range("J58").Value = 1337.1666
$ws.Range("K58").Value = 1296
$ws.Range("L58").Value = 1337.1666
$ws.Range("M58").Value = -1093
$ws.Range("N58").Value = -1743.1666
$ws.Range("H132").Value = 1908.174
$ws.Range("I132").Value = 1922.1818
$ws.Range("K132").Value = 5766.5454
$ws.Range("M132").Value = -3236.5454
$ws.Range("H134").Value = 2169.2632
$ws.Range("I134").Value = 1894.1428
$ws.Range("J134").Value = 2939.6
$ws.Range("K134").Value = 5682.428400000001
$ws.Range("L134").Value = 8818.799999999999
$ws.Range("M134").Value = -3147.428400000001
$ws.Range("N134").Value = -13888.8
$ws.Range("H136").Value = 1305.88
$ws.Range("I136").Value = 1296
$ws.Range("J136").Value = 1337.1666
$ws.Range("K136").Value = 3888
$ws.Range("L136").Value = 4011.4998
$ws.Range("M136").Value = -1338
$ws.Range("N136").Value = -9111.4998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 3950
$ws.Range("I69").Value = 400
$ws.Range("J69").Value = 7500
$ws.Range("K69").Value = 1200
$ws.Range("L69").Value = 22500
$ws.Range("M69").Value = -389
$ws.Range("N69").Value = -24122
$ws.Range("H72").Value = 3950
$ws.Range("I72").Value = 400
$ws.Range("J72").Value = 7500
$ws.Range("K72").Value = 3600
$ws.Range("L72").Value = 67500
$ws.Range("M72").Value = 456
$ws.Range("N72").Value = -75612
$ws.Range("H107").Value = 2194.6667
$ws.Range("I107").Value = 2932.5
$ws.Range("J107").Value = 719
$ws.Range("K107").Value = 8797.5
$ws.Range("L107").Value = 2157
$ws.Range("M107").Value = -6877.5
$ws.Range("N107").Value = -5997
$ws.Range("H113").Value = 1155.0769
$ws.Range("I113").Value = 1060.875
$ws.Range("J113").Value = 1196.9445
$ws.Range("K113").Value = 3182.625
$ws.Range("L113").Value = 3590.8335
$ws.Range("M113").Value = -1012.625
$ws.Range("N113").Value = -7930.833500000001
$ws.Range("H129").Value = 3584.077
$ws.Range("J129").Value = 4848.1113
$ws.Range("L129").Value = 14544.3339
$ws.Range("N129").Value = -24544.3339

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").Value = ""
$ws.Range("H97").Value = 871.8823
$ws.Range("I97").Value = 871.8823
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 871.8823
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -375.8823
$ws.Range("N97").Value = ""
$ws.Range("H102").Value = 2669.2942
$ws.Range("I102").Value = 2311.7144
$ws.Range("K102").Value = 2311.7144
$ws.Range("M102").Value = -689.7143999999998
$ws.Range("H132").Value = 2130
$ws.Range("I132").Value = 1888.3334
$ws.Range("K132").Value = 5665.0002
$ws.Range("M132").Value = -3135.0002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 5755.2354
$ws.Range("I61").Value = 4763.95
$ws.Range("K61").Value = 4763.95
$ws.Range("M61").Value = -4561.95
$ws.Range("H95").Value = 49000
$ws.Range("J95").Value = 49000
$ws.Range("L95").Value = 49000
$ws.Range("N95").Value = -54492
$ws.Range("H100").Value = 8333.333000000001
$ws.Range("I100").Value = 7333.3335
$ws.Range("K100").Value = 7333.3335
$ws.Range("M100").Value = -6792.3335
$ws.Range("H113").Value = 5755.2354
$ws.Range("I113").Value = 4763.95
$ws.Range("K113").Value = 4763.95
$ws.Range("M113").Value = -2593.95
$ws.Range("H132").Value = 4178.875
$ws.Range("J132").Value = 5040
$ws.Range("L132").Value = 15120
$ws.Range("N132").Value = -20180
$ws.Range("H136").Value = 2191.3103
$ws.Range("I136").Value = 2075.55
$ws.Range("J136").Value = 2448.5557
$ws.Range("K136").Value = 6226.650000000001
$ws.Range("L136").Value = 7345.6671
$ws.Range("M136").Value = -3676.650000000001
$ws.Range("N136").Value = -12445.6671

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 25252524
$ws.Range("I5").Value = 100
$ws.Range("J5").Value = 33670000
$ws.Range("K5").Value = 100
$ws.Range("L5").Value = 33670000
$ws.Range("M5").Value = 12
$ws.Range("N5").Value = -33670224
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").Value = ""
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").Value = ""
$ws.Range("H122").Value = 5353.353
$ws.Range("I122").Value = 2786.1428
$ws.Range("J122").Value = 7150.4
$ws.Range("K122").Value = 8358.428400000001
$ws.Range("L122").Value = 21451.2
$ws.Range("M122").Value = -5908.428400000001
$ws.Range("N122").Value = -26351.2
$ws.Range("H132").Value = 6599.75
$ws.Range("I132").Value = 7399.8057
$ws.Range("K132").Value = 22199.4171
$ws.Range("M132").Value = -19669.4171
$ws.Range("H136").Value = 1916.8085
$ws.Range("I136").Value = 1005.06665
$ws.Range("K136").Value = 3015.19995
$ws.Range("M136").Value = -465.1999500000002

Write-Host "applied 260 cell updates"